$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right before the current row 162 (shifts the
# existing rows 162-179 down to 164-181, matching the diff's row renumbering).
$ws.Rows.Item(162).Insert()
$ws.Rows.Item(162).Insert()

# New row 162: Murcott / Especial, Región de O'Higgins
$ws.Cells.Item(162, 1).Value = 4
$ws.Cells.Item(162, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(162, 3).Value = "Los Lagos"
$ws.Cells.Item(162, 4).Value = 44617
$ws.Cells.Item(162, 5).Value = 10
$ws.Cells.Item(162, 6).Value = "Fruta"
$ws.Cells.Item(162, 7).Value = 100102
$ws.Cells.Item(162, 8).Value = "Cítricos"
$ws.Cells.Item(162, 9).Value = 100102004
$ws.Cells.Item(162, 10).Value = "Mandarina"
$ws.Cells.Item(162, 11).Value = "Murcott"
$ws.Cells.Item(162, 12).Value = "Especial"
$ws.Cells.Item(162, 13).Value = 150
$ws.Cells.Item(162, 14).Value = 16000
$ws.Cells.Item(162, 15).Value = 16000
$ws.Cells.Item(162, 16).Value = 16000
$ws.Cells.Item(162, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(162, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(162, 19).Value = 889
$ws.Cells.Item(162, 20).Value = 18

# New row 163: Murcott / Primera, Región de O'Higgins
$ws.Cells.Item(163, 1).Value = 4
$ws.Cells.Item(163, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(163, 3).Value = "Los Lagos"
$ws.Cells.Item(163, 4).Value = 44617
$ws.Cells.Item(163, 5).Value = 10
$ws.Cells.Item(163, 6).Value = "Fruta"
$ws.Cells.Item(163, 7).Value = 100102
$ws.Cells.Item(163, 8).Value = "Cítricos"
$ws.Cells.Item(163, 9).Value = 100102004
$ws.Cells.Item(163, 10).Value = "Mandarina"
$ws.Cells.Item(163, 11).Value = "Murcott"
$ws.Cells.Item(163, 12).Value = "Primera"
$ws.Cells.Item(163, 13).Value = 150
$ws.Cells.Item(163, 14).Value = 15000
$ws.Cells.Item(163, 15).Value = 15000
$ws.Cells.Item(163, 16).Value = 15000
$ws.Cells.Item(163, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(163, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(163, 19).Value = 833
$ws.Cells.Item(163, 20).Value = 18
